$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.76089084148407
$ws.Range("B1").Value = 1.888330578804016
$ws.Range("C1").Value = 2.038364887237549
$ws.Range("D1").Value = 2.919306993484497
$ws.Range("E1").Value = 3.242496728897095
